# Split the "Source PubMed ID (PMID) or doi" column into two separate
# columns - "Source PubMed ID (PMID)" and "Source DOI" - on both the
# "Score Development Samples" and "Evaluation Sample Sets" sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: Score Development Samples ----
$ws5 = $wb.Worksheets.Item("Score Development Samples")

# Insert a new blank column after column O (the old PMID/doi column),
# pushing the old P (Cohort(s)) and Q (Additional Sample/Cohort
# Information) columns one to the right.
$ws5.Columns.Item(16).Insert()

# Re-label the split headers.
$ws5.Range("O1").Value = "Source PubMed ID (PMID)"
$ws5.Range("P1").Value = "Source DOI"

# The PMID values on this sheet were stored as text; store them as
# genuine numbers in their (unchanged) column O.
$ws5.Range("O2").Value = 10000011
$ws5.Range("O3").Value = 10000011
$ws5.Range("O4").Value = 10000012
$ws5.Range("O5").Value = 10000013

# ---- Sheet: Evaluation Sample Sets ----
$ws7 = $wb.Worksheets.Item("Evaluation Sample Sets")

# Same column split here - insert after column O, pushing old P
# (Cohort(s)) and Q (Additional Sample/Cohort Information) to Q/R.
$ws7.Columns.Item(16).Insert()

$ws7.Range("O1").Value = "Source PubMed ID (PMID)"
$ws7.Range("P1").Value = "Source DOI"

# This sheet had no PMIDs, but carries DOIs for each sample set.
$ws7.Range("P2").Value = "10.2021/pgs.1001"
$ws7.Range("P3").Value = "10.2021/pgs.1002"
$ws7.Range("P4").Value = "10.2021/pgs.1003"
